$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 19, shifting existing row 19+ down by one.
$ws.Rows.Item(19).Insert()

# Set the content of the newly inserted row A19.
$ws.Range("A19").Value = '<li><a class="dropdown-item" href="../articles/choropleth_maps.html">Creating choropleth maps</a></li>'
